$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.594.39"
$ws.Range("E2").Value = "  +3.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.373.38"
$ws.Range("E3").Value = "  +4.66%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.93"
$ws.Range("E5").Value = "  +4.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "594.01"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.135"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.422"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.961.66"
$ws.Range("E12").Value = "  +4.70%  "
$ws.Range("E13").Value = "  +1.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.72"
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.615.67"
$ws.Range("E15").Value = "  +3.11%  "
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.376.99"
$ws.Range("E17").Value = "  +5.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "450.72"
$ws.Range("E18").Value = "  +13.63%  "
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.88"
$ws.Range("E20").Value = "  +3.19%  "
$ws.Range("E21").Value = "  +3.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.80"
$ws.Range("E22").Value = "  +5.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000124"
$ws.Range("E24").Value = "  +5.66%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.517.79"
$ws.Range("E25").Value = "  +4.59%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.520"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("E27").Value = "  +4.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.59"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.00"
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.33"
$ws.Range("E31").Value = "  +3.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.67"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.04"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.55"
$ws.Range("E36").Value = "  +5.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "165.04"
$ws.Range("E37").Value = "  +3.03%  "
$ws.Range("E38").Value = "  +3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.27"
$ws.Range("E39").Value = "  +3.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.820"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.57"
$ws.Range("E42").Value = "  +0.74%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.56"
$ws.Range("E43").Value = "  +4.12%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.740.06"
$ws.Range("E44").Value = "  +5.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.83"
$ws.Range("E45").Value = "  +5.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0693"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "344.01"
$ws.Range("E47").Value = "  +3.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.68"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0287"
$ws.Range("E49").Value = "  +3.62%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.10"
$ws.Range("E50").Value = "  +8.54%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.03"
$ws.Range("E51").Value = "  +6.43%  "
